$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group members (new columns B/C added to the header block) ---
$ws.Range("B1").Value = "dragon"
$ws.Range("B2").Value = "vdoo0002"
$ws.Range("C2").Value = "justin-git01"
$ws.Range("B3").Value = "jsan0062"
$ws.Range("C3").Value = "JaySangani"

# --- New meeting row (row 8) replacing the old "Fill in as needed" placeholder ---
$ws.Range("A8").Value = 45189
$ws.Range("A8").NumberFormat = "mm-dd-yy"

$ws.Range("B8").Value = 0.79166666666666663
$ws.Range("B8").NumberFormat = $ws.Range("B7").NumberFormat

$ws.Range("C8").Value = 0.8125
$ws.Range("C8").NumberFormat = $ws.Range("C7").NumberFormat

$ws.Range("D8").Value = "All"

$ws.Range("E8").Value = "Proof-read on all tasks required`nWork allocated for each task`nAgree on deadline for submission on github of task 1`nJustin will fetch the data (Task2 and part of Task3) with help of Jay`nPost on ED if have any question"
$ws.Range("E8").WrapText = $true

$ws.Rows.Item(8).RowHeight = 119
$ws.Columns.Item(4).ColumnWidth = 13.83

# --- Selection moves to C8 ---
$ws.Range("C8").Select()
